# Case_5_11 lines_states.xlsx edit: add two new lines (line7, line8) to the
# "in_service" table, shifting the existing extr1..extr8 rows down by two
# rows and updating their from_bus/to_bus/in_service values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, extend formatting (styles) for the two brand-new rows (16 and 17)
# by copying the format of the last existing data row (15) down into them.
$ws.Range("A15:E15").Copy($ws.Range("A16:E16"))
$ws.Range("A15:E15").Copy($ws.Range("A17:E17"))

# Now fill in the final values. We update from the bottom row upward so that
# each shared string ("extr1".."extr8") stays referenced by some cell at all
# times while it is being "moved" two rows down - this avoids the string
# being dropped and re-added (which would happen if we cleared the old
# cells before writing the new ones).

# Row 17 (was nothing before) <- extr8 values
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false

# Row 16 (was nothing before) <- extr7 values
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

# Row 15 (was extr8) <- extr6 values
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# Row 14 (was extr7) <- extr5 values
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

# Row 13 (was extr6) <- extr4 values
$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# Row 12 (was extr5) <- extr3 values
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $true

# Row 11 (was extr4) <- extr2 values
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# Row 10 (was extr3) <- extr1 values
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# Row 9 (was extr2) <- line8 values
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Row 8 (was extr1) <- line7 values
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true
